$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 11

$ws.Range("D3").Value = 10.35
$ws.Range("E3").Value = 10.49
$ws.Range("F3").Value = 9.57

$ws.Range("C4").Value = 9.65
$ws.Range("E4").Value = 10.57
$ws.Range("F4").Value = 10.3

$ws.Range("B5").Value = 9
$ws.Range("C5").Value = 9.51
$ws.Range("D5").Value = 9.43
$ws.Range("F5").Value = 10.38
$ws.Range("H5").Value = 8.52
$ws.Range("I5").Value = 8

$ws.Range("C6").Value = 10.43
$ws.Range("D6").Value = 9.699999999999999
$ws.Range("E6").Value = 9.619999999999999
$ws.Range("H6").Value = 10.69

$ws.Range("H7").Value = 9.73
$ws.Range("J7").Value = 9.359999999999999

$ws.Range("E8").Value = 11.48
$ws.Range("F8").Value = 9.31
$ws.Range("G8").Value = 10.27

$ws.Range("E9").Value = 12

$ws.Range("G10").Value = 10.64
